# case-control analysis for study periods.
# Rename "Sheet1" to "Dates" and switch the active/selected tab from
# "country-date" to "Dates", moving the selection on the "Dates" sheet
# to B10 (leaving the "country-date" sheet's selection at E85).

$wb = $excel.ActiveWorkbook

$countryDateSheet = $wb.Worksheets.Item("country-date")
$datesSheet = $wb.Worksheets.Item("Sheet1")

# Rename the second sheet.
$datesSheet.Name = "Dates"

# Make sure the selection on the country-date sheet stays put at E85.
$countryDateSheet.Activate()
$countryDateSheet.Range("E85").Select()

# Activate the renamed "Dates" sheet and move its selection to B10 - this
# becomes the active/selected tab for the workbook.
$datesSheet.Activate()
$datesSheet.Range("B10").Select()
